$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("AA4").Value = "*maa://32658 (71.43)"
$ws.Range("O5").Value = "maa://21919 (95.56), maa://21281 (92.31)"
$ws.Range("O6").Value = "maa://31836 (93.33), maa://30381 (100.0)"
$ws.Range("O7").Value = "maa://22750 (96.97)"
$ws.Range("S10").Value = "maa://27395 (97.28), maa://22755 (87.5), **maa://22756 (40.91), ***maa://21737 (10.61)"
$ws.Range("AE10").Value = "*maa://25021 (55.71), *maa://22733 (58.62), maa://22761 (100.0)"
$ws.Range("G13").Value = "*maa://21248 (75.62), **maa://22728 (47.62)"
$ws.Range("AA14").Value = "maa://22764 (96.43)"
$ws.Range("AE15").Value = "maa://21364 (80.61), *maa://22766 (73.0), *maa://36666 (76.67)"
$ws.Range("G17").Value = "maa://22430 (88.57), maa://39599 (81.25)"
$ws.Range("AA18").Value = "maa://24393 (97.14)"
$ws.Range("C20").Value = "maa://21432 (90.7), maa://25198 (93.62), **maa://20795 (49.59), maa://36680 (100.0)"
$ws.Range("K20").Value = "maa://41331 (86.96)"
$ws.Range("K22").Value = "maa://27127 (86.36), *maa://22751 (77.05)"
$ws.Range("K23").Value = "maa://39756 (92.2), maa://39875 (95.12)"
$ws.Range("W23").Value = "*maa://28503 (61.4)"
$ws.Range("S25").Value = "maa://20109 (92.68), maa://22545 (100.0)"
$ws.Range("C26").Value = ""
$ws.Range("W28").Value = "maa://39929 (85.8), ***maa://39723 (15.15), **maa://41749 (50.0)"
$ws.Range("K29").Value = "maa://28432 (93.36), *maa://28440 (72.5), maa://31400 (100.0), *maa://28650 (66.67)"
$ws.Range("S32").Value = "maa://41108 (85.71), maa://41238 (94.12)"
$ws.Range("K35").Value = "maa://41296 (96.55)"
$ws.Range("G39").Value = "maa://25199 (86.11), maa://36670 (88.52), maa://30434 (88.89), ***maa://25036 (16.0)"
$ws.Range("G57").Value = "maa://25176 (97.78)"
